$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") is bumped by one day (46072 -> 46073) for every data row (2..13).
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 3).Value2 = 46073
}

# Rows 6, 9, 10, 11, 12, 13 have their A (Beteckning), B (Datum) and G (Area) values
# re-shuffled among themselves (the underlying records changed order), while all
# other columns for those rows stay the same.
$rows = @(6, 9, 10, 11, 12, 13)

$colA = @{}
$colB = @{}
$colG = @{}
foreach ($r in $rows) {
    $colA[$r] = $ws.Cells.Item($r, 1).Value2
    $colB[$r] = $ws.Cells.Item($r, 2).Value2
    $colG[$r] = $ws.Cells.Item($r, 7).Value2
}

# New order of records (by originating row) placed into each target row.
$mapping = @{
    6  = 13
    9  = 10
    10 = 12
    11 = 9
    12 = 11
    13 = 6
}

foreach ($target in $rows) {
    $source = $mapping[$target]
    $ws.Cells.Item($target, 1).Value2 = $colA[$source]
    $ws.Cells.Item($target, 2).Value2 = $colB[$source]
    $ws.Cells.Item($target, 7).Value2 = $colG[$source]
}
